# Update "想去人数" (F column) figures that changed between scrapes.
#
# Sheet 1 "展览"   (rows 4-13)
# Sheet 2 "演出"   (row 2)
# Sheet 4 "全部类型" (rows 4-13, and row 14 which mirrors 演出's row 2)
#
# Sheet 3 "本地生活" has no data rows and needs no change.

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    4  = 63
    5  = 511
    6  = 6810
    7  = 190
    8  = 151
    9  = 1038
    10 = 387
    11 = 130
    12 = 187
    13 = 564
}

# Sheet 1: 展览
$wsExhibition = $wb.Worksheets.Item(1)
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F" + $row).Value = $exhibitionUpdates[$row]
}

# Sheet 2: 演出
$wsPerformance = $wb.Worksheets.Item(2)
$wsPerformance.Range("F2").Value = 19

# Sheet 4: 全部类型 (combined sheet — same exhibition rows, plus the
# performance row appended at the end as row 14)
$wsAll = $wb.Worksheets.Item(4)
foreach ($row in $exhibitionUpdates.Keys) {
    $wsAll.Range("F" + $row).Value = $exhibitionUpdates[$row]
}
$wsAll.Range("F14").Value = 19
